$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells keep their original text formatting
# (prevents Excel from auto-converting numeric-looking strings to
# floating point numbers and losing precision / trailing zeros).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '63.953.73'
$ws.Range('E2').Value = '  +3.48%  '
$ws.Range('D3').Value = '3.052.98'
$ws.Range('E3').Value = '  +2.46%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '556.43'
$ws.Range('E5').Value = '  +2.78%  '
$ws.Range('D6').Value = '142.46'
$ws.Range('E6').Value = '  +4.92%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').Value = '3.049.35'
$ws.Range('E8').Value = '  +2.54%  '
$ws.Range('D9').Value = '0.512'
$ws.Range('E9').Value = '  +5.42%  '
$ws.Range('E10').Value = '  +6.26%  '
$ws.Range('E11').Value = '  -9.57%  '
$ws.Range('E12').Value = '  +8.12%  '
$ws.Range('E13').Value = '  +6.08%  '
$ws.Range('D14').Value = '35.08'
$ws.Range('E14').Value = '  +4.12%  '
$ws.Range('D15').Value = '3.546.85'
$ws.Range('E15').Value = '  +2.94%  '
$ws.Range('D16').Value = '63.995.91'
$ws.Range('E16').Value = '  +3.53%  '
$ws.Range('D17').Value = '3.052.53'
$ws.Range('E17').Value = '  +2.51%  '
$ws.Range('E18').Value = '  +1.94%  '
$ws.Range('E19').Value = '  +3.49%  '
$ws.Range('D20').Value = '477.05'
$ws.Range('E20').Value = '  +2.70%  '
$ws.Range('D21').Value = '13.99'
$ws.Range('E21').Value = '  +5.28%  '
$ws.Range('D22').Value = '0.678'
$ws.Range('E22').Value = '  +4.43%  '
$ws.Range('D23').Value = '7.58'
$ws.Range('E23').Value = '  +6.69%  '
$ws.Range('D24').Value = '14.22'
$ws.Range('E24').Value = '  +14.34%  '
$ws.Range('D25').Value = '81.85'
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('E27').Value = '  +2.81%  '
$ws.Range('E28').Value = '  +5.29%  '
$ws.Range('E29').Value = '  +2.49%  '
$ws.Range('E30').Value = '  +0.22%  '
$ws.Range('E31').Value = '  +4.16%  '
$ws.Range('E32').Value = '  +2.19%  '
$ws.Range('E33').Value = '  +5.23%  '
$ws.Range('D34').Value = '5.63'
$ws.Range('E34').Value = '  +3.21%  '
$ws.Range('D35').Value = '6.19'
$ws.Range('E35').Value = '  +7.10%  '
$ws.Range('D36').Value = '54.72'
$ws.Range('E36').Value = '  +1.84%  '
$ws.Range('D37').Value = '0.0406'
$ws.Range('E37').Value = '  +5.72%  '
$ws.Range('D38').Value = '442.17'
$ws.Range('E38').Value = '  -1.47%  '
$ws.Range('D39').Value = '0.0806'
$ws.Range('E39').Value = '  +0.55%  '
$ws.Range('E40').Value = '  +15.01%  '
$ws.Range('D41').Value = '2.972.63'
$ws.Range('E41').Value = '  +1.18%  '
$ws.Range('E42').Value = '  +3.15%  '
$ws.Range('D44').Value = '27.65'
$ws.Range('E44').Value = '  +3.70%  '
$ws.Range('D45').Value = '0.260'
$ws.Range('E45').Value = '  +6.38%  '
$ws.Range('E46').Value = '  +9.24%  '
$ws.Range('E47').Value = '  -0.02%  '
$ws.Range('E48').Value = '  +4.93%  '
$ws.Range('D49').Value = '0.0₃0514'
$ws.Range('E49').Value = '  +5.74%  '
$ws.Range('D50').Value = '117.11'
$ws.Range('E50').Value = '  +2.77%  '
$ws.Range('D51').Value = '2.06'
$ws.Range('E51').Value = '  +3.74%  '
